$wb = $excel.ActiveWorkbook

# --- Sheet1: three_line ---
$ws1 = $wb.Worksheets.Item("three_line")
$ws1.Range("A194").Value = 45441.55208333334
$ws1.Range("A194").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B194").Value = "11-06-2024 10:15:00"
$ws1.Range("C194").Value = "hour"
$ws1.Range("D194").Value = "TATASTEEL.NS"
$ws1.Range("E194").Value = 45415.38541666666
$ws1.Range("E194").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("F194").Value = 170.75
$ws1.Range("G194").Value = 45434.38541666666
$ws1.Range("G194").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("H194").Value = 175.4499969482422
$ws1.Range("I194").Value = 45440.38541666666
$ws1.Range("I194").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("J194").Value = 177.5
$ws1.Range("K194").Value = "High"
$ws1.Range("L194").Value = "11/06/2024 04:47:04"
$ws1.Range("A195").Value = 45441.55208333334
$ws1.Range("A195").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B195").Value = "11-06-2024 10:15:00"
$ws1.Range("C195").Value = "hour"
$ws1.Range("D195").Value = "TATASTEEL.NS"
$ws1.Range("E195").Value = 45415.38541666666
$ws1.Range("E195").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("F195").Value = 170.75
$ws1.Range("G195").Value = 45436.38541666666
$ws1.Range("G195").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("H195").Value = 177.5500030517578
$ws1.Range("I195").Value = 45440.38541666666
$ws1.Range("I195").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("J195").Value = 177.5
$ws1.Range("K195").Value = "High"
$ws1.Range("L195").Value = "11/06/2024 04:47:04"
$ws1.Range("A196").Value = 45447.55208333334
$ws1.Range("A196").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B196").Value = "11-06-2024 09:15:00"
$ws1.Range("C196").Value = "hour"
$ws1.Range("D196").Value = "DHAMPURSUG.NS"
$ws1.Range("E196").Value = 45436.38541666666
$ws1.Range("E196").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("F196").Value = 225
$ws1.Range("G196").Value = 45439.38541666666
$ws1.Range("G196").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("H196").Value = 225
$ws1.Range("I196").Value = 45446.38541666666
$ws1.Range("I196").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("J196").Value = 225.25
$ws1.Range("K196").Value = "High"
$ws1.Range("L196").Value = "11/06/2024 04:47:04"
$ws1.Range("A197").Value = 45419.55208333334
$ws1.Range("A197").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B197").Value = "11-06-2024 09:15:00"
$ws1.Range("C197").Value = "hour"
$ws1.Range("D197").Value = "JKIL.NS"
$ws1.Range("E197").Value = 45394.38541666666
$ws1.Range("E197").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("F197").Value = 644.5999755859375
$ws1.Range("G197").Value = 45411.59375
$ws1.Range("G197").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("H197").Value = 690.2999877929688
$ws1.Range("I197").Value = 45418.38541666666
$ws1.Range("I197").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("J197").Value = 704.75
$ws1.Range("K197").Value = "High"
$ws1.Range("L197").Value = "11/06/2024 04:47:04"
$ws1.Range("A198").Value = 45408.38541666666
$ws1.Range("A198").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B198").Value = "11-06-2024 09:15:00"
$ws1.Range("C198").Value = "hour"
$ws1.Range("D198").Value = "ANDHRAPET.BO"
$ws1.Range("E198").Value = 45377.38541666666
$ws1.Range("E198").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("F198").Value = 87.5
$ws1.Range("G198").Value = 45386.38541666666
$ws1.Range("G198").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("H198").Value = 91.94999694824219
$ws1.Range("I198").Value = 45394.38541666666
$ws1.Range("I198").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("J198").Value = 94.9000015258789
$ws1.Range("K198").Value = "High"
$ws1.Range("L198").Value = "11/06/2024 04:47:04"

# --- Sheet2: two_line ---
$ws2 = $wb.Worksheets.Item("two_line")
$ws2.Range("A36").Value = 45447.55208333334
$ws2.Range("A36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B36").Value = "11-06-2024 10:15:00"
$ws2.Range("C36").Value = "hour"
$ws2.Range("D36").Value = "KOTAKBANK.NS"
$ws2.Range("E36").Value = 45433.38541666666
$ws2.Range("E36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("F36").Value = 1714.900024414062
$ws2.Range("G36").Value = 45446.38541666666
$ws2.Range("G36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("H36").Value = 1726.449951171875
$ws2.Range("I36").Value = "High"
$ws2.Range("J36").Value = "11/06/2024 04:47:04"
$ws2.Range("A37").Value = 45439.55208333334
$ws2.Range("A37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B37").Value = "11-06-2024 10:15:00"
$ws2.Range("C37").Value = "hour"
$ws2.Range("D37").Value = "DLF.NS"
$ws2.Range("E37").Value = 45429.46875
$ws2.Range("E37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("F37").Value = 859.5
$ws2.Range("G37").Value = 45434.38541666666
$ws2.Range("G37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("H37").Value = 860
$ws2.Range("I37").Value = "High"
$ws2.Range("J37").Value = "11/06/2024 04:47:04"

# --- Sheet3: ph_pl_breakout_line ---
$ws3 = $wb.Worksheets.Item("ph_pl_breakout_line")
$ws3.Range("A684").Value = "HERCULES.NS"
$ws3.Range("B684").Value = 45450.42708333334
$ws3.Range("B684").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C684").Value = 493.3999938964844
$ws3.Range("D684").Value = 463.75
$ws3.Range("E684").Value = 480.2000122070312
$ws3.Range("F684").Value = "High"
$ws3.Range("G684").Value = 493.3999938964844
$ws3.Range("H684").Value = "hour"
$ws3.Range("I684").Value = "11-06-2024 09:15:00"
$ws3.Range("J684").Value = 497
$ws3.Range("K684").Value = 491.5
$ws3.Range("L684").Value = "11/06/2024 04:47:04"
$ws3.Range("A685").Value = "ARTSONEN.BO"
$ws3.Range("B685").Value = 45439.46875
$ws3.Range("B685").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C685").Value = 189.3999938964844
$ws3.Range("D685").Value = 189.3000030517578
$ws3.Range("E685").Value = 189.3999938964844
$ws3.Range("F685").Value = "High"
$ws3.Range("G685").Value = 189.3999938964844
$ws3.Range("H685").Value = "hour"
$ws3.Range("I685").Value = "11-06-2024 09:15:00"
$ws3.Range("J685").Value = 193.3000030517578
$ws3.Range("K685").Value = 184.6000061035156
$ws3.Range("L685").Value = "11/06/2024 04:47:04"
$ws3.Range("A686").Value = "ARTSONEN.BO"
$ws3.Range("B686").Value = 45439.51041666666
$ws3.Range("B686").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C686").Value = 189.3999938964844
$ws3.Range("D686").Value = 185.4499969482422
$ws3.Range("E686").Value = 185.4499969482422
$ws3.Range("F686").Value = "High"
$ws3.Range("G686").Value = 189.3999938964844
$ws3.Range("H686").Value = "hour"
$ws3.Range("I686").Value = "11-06-2024 09:15:00"
$ws3.Range("J686").Value = 193.3000030517578
$ws3.Range("K686").Value = 184.6000061035156
$ws3.Range("L686").Value = "11/06/2024 04:47:04"
$ws3.Range("A687").Value = "ARTSONEN.BO"
$ws3.Range("B687").Value = 45439.55208333334
$ws3.Range("B687").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C687").Value = 189.3999938964844
$ws3.Range("D687").Value = 185
$ws3.Range("E687").Value = 188
$ws3.Range("F687").Value = "High"
$ws3.Range("G687").Value = 189.3999938964844
$ws3.Range("H687").Value = "hour"
$ws3.Range("I687").Value = "11-06-2024 09:15:00"
$ws3.Range("J687").Value = 193.3000030517578
$ws3.Range("K687").Value = 184.6000061035156
$ws3.Range("L687").Value = "11/06/2024 04:47:04"
$ws3.Range("A688").Value = "TRIL.BO"
$ws3.Range("B688").Value = 45449.38541666666
$ws3.Range("B688").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C688").Value = 750
$ws3.Range("D688").Value = 700
$ws3.Range("E688").Value = 722
$ws3.Range("F688").Value = "High"
$ws3.Range("G688").Value = 750
$ws3.Range("H688").Value = "hour"
$ws3.Range("I688").Value = "11-06-2024 09:15:00"
$ws3.Range("J688").Value = 769.1500244140625
$ws3.Range("K688").Value = 735
$ws3.Range("L688").Value = "11/06/2024 04:47:04"
$ws3.Range("A689").Value = "TRIL.BO"
$ws3.Range("B689").Value = 45449.42708333334
$ws3.Range("B689").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C689").Value = 750
$ws3.Range("D689").Value = 728
$ws3.Range("E689").Value = 740
$ws3.Range("F689").Value = "High"
$ws3.Range("G689").Value = 750
$ws3.Range("H689").Value = "hour"
$ws3.Range("I689").Value = "11-06-2024 09:15:00"
$ws3.Range("J689").Value = 769.1500244140625
$ws3.Range("K689").Value = 735
$ws3.Range("L689").Value = "11/06/2024 04:47:04"
$ws3.Range("A690").Value = "ASMTEC.BO"
$ws3.Range("B690").Value = 45446.38541666666
$ws3.Range("B690").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C690").Value = 1100.400024414062
$ws3.Range("D690").Value = 1075
$ws3.Range("E690").Value = 1100.400024414062
$ws3.Range("F690").Value = "High"
$ws3.Range("G690").Value = 1100.400024414062
$ws3.Range("H690").Value = "hour"
$ws3.Range("I690").Value = "11-06-2024 09:15:00"
$ws3.Range("J690").Value = 1141.75
$ws3.Range("K690").Value = 1087.400024414062
$ws3.Range("L690").Value = "11/06/2024 04:47:04"
$ws3.Range("A691").Value = "SPELS.BO"
$ws3.Range("B691").Value = 45434.42708333334
$ws3.Range("B691").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C691").Value = 141.6999969482422
$ws3.Range("D691").Value = 141.6999969482422
$ws3.Range("E691").Value = 141.6999969482422
$ws3.Range("F691").Value = "High"
$ws3.Range("G691").Value = 141.6999969482422
$ws3.Range("H691").Value = "hour"
$ws3.Range("I691").Value = "11-06-2024 09:15:00"
$ws3.Range("J691").Value = 143.8000030517578
$ws3.Range("K691").Value = 139.6499938964844
$ws3.Range("L691").Value = "11/06/2024 04:47:04"
$ws3.Range("A692").Value = "SPELS.BO"
$ws3.Range("B692").Value = 45434.46875
$ws3.Range("B692").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C692").Value = 141.6999969482422
$ws3.Range("D692").Value = 141.6999969482422
$ws3.Range("E692").Value = 141.6999969482422
$ws3.Range("F692").Value = "High"
$ws3.Range("G692").Value = 141.6999969482422
$ws3.Range("H692").Value = "hour"
$ws3.Range("I692").Value = "11-06-2024 09:15:00"
$ws3.Range("J692").Value = 143.8000030517578
$ws3.Range("K692").Value = 139.6499938964844
$ws3.Range("L692").Value = "11/06/2024 04:47:04"
$ws3.Range("A693").Value = "CENTRALBK.BO"
$ws3.Range("B693").Value = 45440.51041666666
$ws3.Range("B693").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C693").Value = 65.98999786376953
$ws3.Range("D693").Value = 63.90000152587891
$ws3.Range("E693").Value = 65.37999725341797
$ws3.Range("F693").Value = "Low"
$ws3.Range("G693").Value = 63.90000152587891
$ws3.Range("H693").Value = "hour"
$ws3.Range("I693").Value = "11-06-2024 09:15:00"
$ws3.Range("J693").Value = 63.47999954223633
$ws3.Range("K693").Value = 64
$ws3.Range("L693").Value = "11/06/2024 04:47:04"
$ws3.Range("A694").Value = "RPEL.BO"
$ws3.Range("B694").Value = 45439.38541666666
$ws3.Range("B694").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C694").Value = 695.9000244140625
$ws3.Range("D694").Value = 672
$ws3.Range("E694").Value = 686
$ws3.Range("F694").Value = "Low"
$ws3.Range("G694").Value = 672
$ws3.Range("H694").Value = "hour"
$ws3.Range("I694").Value = "11-06-2024 09:15:00"
$ws3.Range("J694").Value = 671.1500244140625
$ws3.Range("K694").Value = 672.0999755859375
$ws3.Range("L694").Value = "11/06/2024 04:47:04"
$ws3.Range("A695").Value = "UNIAUTO.BO"
$ws3.Range("B695").Value = 45433.38541666666
$ws3.Range("B695").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("C695").Value = 168.3999938964844
$ws3.Range("D695").Value = 163
$ws3.Range("E695").Value = 163
$ws3.Range("F695").Value = "High"
$ws3.Range("G695").Value = 168.3999938964844
$ws3.Range("H695").Value = "hour"
$ws3.Range("I695").Value = "11-06-2024 09:15:00"
$ws3.Range("J695").Value = 171.9499969482422
$ws3.Range("K695").Value = 167.8500061035156
$ws3.Range("L695").Value = "11/06/2024 04:47:04"

Write-Host "Edit complete"